$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows: ManufacturerId changes for rows 2 and 6 (B2, B6: 3 -> 8)
$ws.Range("B2").Value = 8
$ws.Range("B6").Value = 8

# New data rows 7-13
$data = @(
    @(6, 8, "gtx1060", "2gb", 350),
    @(7, 8, "gtx980",  "4gb", 300),
    @(8, 8, "gtx970",  "2gb", 250),
    @(9, 6, "r9290",   "4gb", 300),
    @(10, 6, "r9280",  "2gb", 250),
    @(11, 6, "r9390",  "4gb", 350),
    @(12, 6, "r9380",  "2gb", 300)
)

$row = 7
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $ws.Cells.Item($row, 5).Value = $item[4]
    $row++
}

# Update selection to match final state
$ws.Range("E13").Select()
